# RoboRIO Ports and other controls.xlsx
# "Added overrides and finished scaling code"
#
# The "Scaler Left Lift" / "Scaler Right Lift" entries in the RefNum column
# are renamed to the finished "Scaler Lift 1" / "Scaler Lift 2" naming, and
# the active selection is left on the cell the author was last working in.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RoboRIO Ports")

$ws.Range("C6").Value = "Scaler Lift 1"
$ws.Range("C7").Value = "Scaler Lift 2"

$ws.Range("D7").Select() | Out-Null
